# Builder AV Example update:
# - Row 2 keeps its row reference but the underlying trial data changes.
# - Row 3 gets new trial data (its T_SND cell format is normalised to match
#   the rest of the column - it previously used a slightly different font).
# - Rows 4 and 5 are updated to a new repeating trial pattern.
# - Two additional trial rows (6 and 7) are appended, continuing that
#   pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 ----------------------------------------------------------------
$ws.Range("C2").Value = "2.0"
$ws.Range("D2").Value = "go.wav"
$ws.Range("E2").Value = "."
$ws.Range("F2").Value = "Line768.png"

# ---- Row 3 ------------------------------------------------------------
$ws.Range("C3").Value = "2.0"
# D3 used an outlier font/style in the original sheet; align it with the
# rest of the (text-formatted) column before writing the new value.
$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").Value = "a2z_upper_inst.wav"
$ws.Range("E3").Value = "."
$ws.Range("F3").Value = "Line768x3.png"

# ---- Row 4 ------------------------------------------------------------
$ws.Range("C4").Value = "1"
$ws.Range("C4").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D4").Value = "drawcircles.wav"
$ws.Range("E4").Value = "CirclesCCW.mp4"
$ws.Range("F4").Value = "green_dot.png"

# ---- Row 5 ------------------------------------------------------------
$ws.Range("C5").Value = "1"
$ws.Range("C5").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("D5").Value = "drawcircles.wav"
$ws.Range("E5").Value = "CirclesCCW.mp4"
$ws.Range("F5").Value = "green_dot.png"

# ---- Row 6 (new, same pattern as rows 4/7 - taller custom row height) ----
$ws.Range("A5:H5").Copy()
$ws.Range("A6:H6").PasteSpecial(-4122)
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = -1
$ws.Range("C6").Value = "1"
$ws.Range("D6").Value = "drawcircles.wav"
$ws.Range("E6").Value = "CirclesCCW.mp4"
$ws.Range("F6").Value = "green_dot.png"
$ws.Range("G6").Value = -1.1
$ws.Range("H6").Value = -1.1
$ws.Rows.Item(6).RowHeight = 21.75

# ---- Row 7 (new) ------------------------------------------------------
$ws.Range("A5:H5").Copy()
$ws.Range("A7:H7").PasteSpecial(-4122)
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = -1
$ws.Range("C7").Value = "1"
$ws.Range("D7").Value = "drawcircles.wav"
$ws.Range("E7").Value = "CirclesCCW.mp4"
$ws.Range("F7").Value = "green_dot.png"
$ws.Range("G7").Value = -1.1
$ws.Range("H7").Value = -1.1

# ---- Selection matches the saved view in the updated workbook -------------
$ws.Range("G12").Select()
